$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift all timestamps in column A (rows 2-97) forward by 4 days
for ($r = 2; $r -le 97; $r++) {
    $cell = $ws.Cells.Item($r, 1)
    $cell.Value = $cell.Value2 + 4
}

# Update forecasted solar production values in column B (rows 21-41)
$ws.Cells.Item(21, 2).Value = 8
$ws.Cells.Item(22, 2).Value = 28
$ws.Cells.Item(23, 2).Value = 61
$ws.Cells.Item(24, 2).Value = 99
$ws.Cells.Item(25, 2).Value = 145
$ws.Cells.Item(26, 2).Value = 229
$ws.Cells.Item(27, 2).Value = 316
$ws.Cells.Item(28, 2).Value = 389
$ws.Cells.Item(29, 2).Value = 470
$ws.Cells.Item(30, 2).Value = 599
$ws.Cells.Item(31, 2).Value = 729
$ws.Cells.Item(32, 2).Value = 828
$ws.Cells.Item(33, 2).Value = 898
$ws.Cells.Item(34, 2).Value = 1021
$ws.Cells.Item(35, 2).Value = 1138
$ws.Cells.Item(36, 2).Value = 1250
$ws.Cells.Item(37, 2).Value = 1324
$ws.Cells.Item(38, 2).Value = 1428
$ws.Cells.Item(39, 2).Value = 1498
$ws.Cells.Item(40, 2).Value = 1540
$ws.Cells.Item(41, 2).Value = 1560
